$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = $null

$ws.Range("H15").Value = 2367.2273
$ws.Range("I15").Value = 2367.2273
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 7101.6819
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -6932.6819

$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").Value = $null

$ws.Range("H113").Value = 6034.1816
$ws.Range("I113").Value = 6374.625
$ws.Range("J113").Value = 5126.3335
$ws.Range("K113").Value = 6374.625
$ws.Range("L113").Value = 5126.3335
$ws.Range("M113").Value = -3120.625
$ws.Range("N113").Value = -11634.3335

$ws.Range("H137").Value = 1864.68
$ws.Range("I137").Value = 1857.6818
$ws.Range("J137").Value = 1916
$ws.Range("K137").Value = 5573.0454
$ws.Range("L137").Value = 5748
$ws.Range("M137").Value = -3023.0454
$ws.Range("N137").Value = -10848

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 5052.5
$ws.Range("I3").Value = 105
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 105
$ws.Range("L3").Value = 10000
$ws.Range("M3").Value = 10
$ws.Range("N3").Value = -10230

$ws.Range("H32").Value = 5700.643
$ws.Range("I32").Value = 5700.643
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 5700.643
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -5413.643

$ws.Range("H102").Value = 1119.4375
$ws.Range("I102").Value = 1283.3846
$ws.Range("J102").Value = 409
$ws.Range("K102").Value = 1283.3846
$ws.Range("L102").Value = 409
$ws.Range("M102").Value = 338.6153999999999
$ws.Range("N102").Value = -3653

$ws.Range("H110").Value = 3118.5
$ws.Range("I110").Value = 1513.3334
$ws.Range("J110").Value = 12749.5
$ws.Range("K110").Value = 1513.3334
$ws.Range("L110").Value = 12749.5
$ws.Range("M110").Value = 531.6666
$ws.Range("N110").Value = -16839.5

$ws.Range("H132").Value = 3369.3
$ws.Range("I132").Value = 2670.8572
$ws.Range("J132").Value = 4999
$ws.Range("K132").Value = 8012.571599999999
$ws.Range("L132").Value = 14997
$ws.Range("M132").Value = -5482.571599999999
$ws.Range("N132").Value = -20057

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 642.5714
$ws.Range("I107").Value = 619.6
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 619.6
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1300.4
$ws.Range("N107").Value = -4540

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1680
$ws.Range("I16").Value = 1680
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 1680
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -1393
$ws.Range("N16").Value = $null

$ws.Range("H31").Value = 1649.5834
$ws.Range("I31").Value = 1610.6666
$ws.Range("J31").Value = 1766.3334
$ws.Range("K31").Value = 1610.6666
$ws.Range("L31").Value = 1766.3334
$ws.Range("M31").Value = -1315.6666
$ws.Range("N31").Value = -2356.3334

$ws.Range("H34").Value = 1649.5834
$ws.Range("I34").Value = 1610.6666
$ws.Range("J34").Value = 1766.3334
$ws.Range("K34").Value = 1610.6666
$ws.Range("L34").Value = 1766.3334
$ws.Range("M34").Value = -1408.6666
$ws.Range("N34").Value = -2170.3334

$ws.Range("H107").Value = 782.38464
$ws.Range("I107").Value = 789.25
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 789.25
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1130.75
$ws.Range("N107").Value = -4540

$ws.Range("H113").Value = 1680
$ws.Range("I113").Value = 1680
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1680
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 490
$ws.Range("N113").Value = $null

$ws.Range("H132").Value = 3662.8333
$ws.Range("I132").Value = 3239.889
$ws.Range("J132").Value = 4931.6665
$ws.Range("K132").Value = 9719.667000000001
$ws.Range("L132").Value = 14794.9995
$ws.Range("M132").Value = -7189.667000000001
$ws.Range("N132").Value = -19854.9995

$ws.Range("H134").Value = 10488
$ws.Range("I134").Value = 10737.25
$ws.Range("J134").Value = 9989.5
$ws.Range("K134").Value = 32211.75
$ws.Range("L134").Value = 29968.5
$ws.Range("M134").Value = -29676.75
$ws.Range("N134").Value = -35038.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = $null

$ws.Range("H121").Value = 323.42856
$ws.Range("I121").Value = 382.8
$ws.Range("J121").Value = 175
$ws.Range("K121").Value = 1148.4
$ws.Range("L121").Value = 525
$ws.Range("M121").Value = 161.5999999999999
$ws.Range("N121").Value = -3145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 410.66666
$ws.Range("I2").Value = 461.2857
$ws.Range("J2").Value = 339.8
$ws.Range("K2").Value = 461.2857
$ws.Range("L2").Value = 339.8
$ws.Range("M2").Value = -348.2857
$ws.Range("N2").Value = -565.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1507
$ws.Range("I7").Value = 1507
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 1507
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -1395
$ws.Range("N7").Value = $null

$ws.Range("H20").Value = 5000000
$ws.Range("I20").Value = 5000000
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 5000000
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = -4999774

$ws.Range("H76").Value = 74880
$ws.Range("I76").Value = 0
$ws.Range("J76").Value = 74880
$ws.Range("K76").Value = 0
$ws.Range("L76").Value = 74880
$ws.Range("N76").Value = -75556

$ws.Range("H79").Value = 74880
$ws.Range("I79").Value = 0
$ws.Range("J79").Value = 74880
$ws.Range("K79").Value = 0
$ws.Range("L79").Value = 74880
$ws.Range("N79").Value = -77220

$ws.Range("H93").Value = 2437.8572
$ws.Range("I93").Value = 1701.6666
$ws.Range("J93").Value = 2990
$ws.Range("K93").Value = 1701.6666
$ws.Range("L93").Value = 2990
$ws.Range("M93").Value = -453.6666
$ws.Range("N93").Value = -5486

$ws.Range("H126").Value = 1507
$ws.Range("I126").Value = 1507
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4521
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -2051
$ws.Range("N126").Value = $null

$ws.Range("H132").Value = 3462.6428
$ws.Range("I132").Value = 2458
$ws.Range("J132").Value = 5974.25
$ws.Range("K132").Value = 7374
$ws.Range("L132").Value = 17922.75
$ws.Range("M132").Value = -4844
$ws.Range("N132").Value = -22982.75

$ws.Range("H136").Value = 3332.1667
$ws.Range("I136").Value = 3332.1667
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 9996.500100000001
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7446.500100000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 10000
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10346

$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").Value = $null

$ws.Range("H81").Value = 4000.4
$ws.Range("I81").Value = 1667.1111
$ws.Range("J81").Value = 25000
$ws.Range("K81").Value = 3334.2222
$ws.Range("L81").Value = 50000
$ws.Range("M81").Value = -2273.2222
$ws.Range("N81").Value = -52122

$ws.Range("H84").Value = 4000.4
$ws.Range("I84").Value = 1667.1111
$ws.Range("J84").Value = 25000
$ws.Range("K84").Value = 16671.111
$ws.Range("L84").Value = 250000
$ws.Range("M84").Value = -11367.111
$ws.Range("N84").Value = -260608
